$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("boletim_incompleto")

# Replace the text (comma-decimal) grades with real numeric values,
# and fix D3 which actually changes from 6,9 -> 6.7
$ws.Range("B2").Value = 5.6
$ws.Range("C2").Value = 6.7
$ws.Range("D2").Value = 6.7
$ws.Range("E2").Value = 5.8

$ws.Range("B3").Value = 7.2
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 6.7
$ws.Range("E3").Value = 9

$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 5.9
$ws.Range("D4").Value = 8.7
$ws.Range("E4").Value = 6.9

$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 5.6
$ws.Range("E5").Value = 7

$ws.Range("B6").Value = 8.6
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 7.9
$ws.Range("E6").Value = 8

$ws.Range("B7").Value = 9.5
$ws.Range("C7").Value = 8.5
$ws.Range("D7").Value = 9
$ws.Range("E7").Value = 7

# Remove the stray empty/formatted cell that was left at row 11
# (clearing it fully empties row 11, so it drops out of the sheet
# and the used range shrinks back down to row 7)
$ws.Range("B11").Clear()

# Select E5 to match the final cursor location, applying the "typed" style
$ws.Range("E5").Select()
